# Fruta / hortaliza, semanal
# Updates a batch of weekly Espárragos records: dates, volumes, min/max/avg
# prices and origin for several rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @{
    2  = @{ D = 44167; J = 1430; K = 1200; M = 1248; P = 1248 }
    3  = @{ D = 44167; J = 350 }
    6  = @{ D = 44168; J = 1200 }
    7  = @{ D = 44168 }
    8  = @{ D = 44176; J = 2500; K = 1200; M = 1256; P = 1256 }
    9  = @{ D = 44176; J = 1500 }
    10 = @{ D = 44161; J = 1600 }
    11 = @{ D = 44161; J = 1850 }
    16 = @{ D = 44169; J = 950 }
    17 = @{ D = 44169; J = 800 }
    18 = @{ D = 44172; J = 600; K = 1300; L = 1300; M = 1300; P = 1300 }
    19 = @{ D = 44172; J = 550 }
    20 = @{ D = 44160; J = 750; K = 1300; M = 1300; P = 1300 }
    21 = @{ D = 44160; J = 850 }
    22 = @{ D = 44159; J = 1100; K = 1300; M = 1300; P = 1300 }
    23 = @{ D = 44159; J = 800 }
    24 = @{ D = 44162; J = 1200; K = 1300; L = 1300; M = 1300; O = "Provincia de Quillota"; P = 1300 }
    25 = @{ I = "Segunda"; J = 800; K = 1000; L = 1000; M = 1000; P = 1000 }
    26 = @{ D = 44165; I = "Primera"; J = 720; K = 1200; L = 1200; M = 1200; P = 1200 }
    27 = @{ D = 44165; I = "Segunda"; J = 750; K = 1000; L = 1000; M = 1000; P = 1000 }
    28 = @{ D = 44179; I = "Primera"; J = 980; K = 1200; L = 1200; M = 1200; O = "Región Metropolitana"; P = 1200 }
}

foreach ($r in $rowUpdates.Keys) {
    $cols = $rowUpdates[$r]
    foreach ($c in $cols.Keys) {
        $addr = "$c$r"
        $ws.Range($addr).Value = $cols[$c]
    }
}
